$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72, shifting existing rows 72-157 down to 73-158
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new data point
$ws.Cells.Item(72, 1).Value = 11
$ws.Cells.Item(72, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(72, 3).Value = "Bíobío"
$ws.Cells.Item(72, 4).Value = 44539
$ws.Cells.Item(72, 4).Style = $ws.Cells.Item(73, 4).Style
$ws.Cells.Item(72, 4).NumberFormat = $ws.Cells.Item(73, 4).NumberFormat
$ws.Cells.Item(72, 5).Value = 8
$ws.Cells.Item(72, 6).Value = 100114013
$ws.Cells.Item(72, 7).Value = "Zanahoria"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 250
$ws.Cells.Item(72, 11).Value = 5500
$ws.Cells.Item(72, 12).Value = 6000
$ws.Cells.Item(72, 13).Value = 5700
$ws.Cells.Item(72, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(72, 15).Value = "Chillán"
$ws.Cells.Item(72, 16).Value = 285
$ws.Cells.Item(72, 17).Value = 20
$ws.Cells.Item(72, 18).Value = "Hortaliza"
